$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Grab a clean copy (as WordOpenXML) of the "Meta description" run
#    structure so we can re-use its run layout (leading empty run +
#    a bold run) for the new heading paragraph that needs to be
#    inserted just before the final "Create a feature image..."
#    paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaXml = $metaPara.Range.WordOpenXML

# Strip the revision-save-id / rsid attributes that Word stamps onto
# the paragraph when round-tripping through WordOpenXML - the rest of
# this document doesn't carry them, so keep the XML consistent.
$cleaned = $metaXml -replace ' w14:paraId="[0-9A-Fa-f]+"', ''
$cleaned = $cleaned -replace ' w14:textId="[0-9A-Fa-f]+"', ''
$cleaned = $cleaned -replace ' w:rsidR="[0-9A-Fa-f]+"', ''
$cleaned = $cleaned -replace ' w:rsidRDefault="[0-9A-Fa-f]+"', ''

# ------------------------------------------------------------------
# 2) Insert that paragraph markup right before the last paragraph in
#    the document (the "Create a feature image..." paragraph).
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertPoint.InsertXML($cleaned)

# ------------------------------------------------------------------
# 3) Remove the original "Meta description" paragraph (still the 2nd
#    paragraph of the document at this point).
# ------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Delete()

# InsertXML also produces a trailing empty paragraph (the end of the
# inserted "story"); remove it too. After step 3 it now sits right
# before the "Create a feature image..." paragraph.
$countAfterInsert = $d.Paragraphs.Count
$emptyPara = $d.Paragraphs.Item($countAfterInsert - 1)
$emptyPara.Range.Delete()

# ------------------------------------------------------------------
# 4) Turn the freshly-inserted paragraph's text from
#       "Meta description: Immersive gameplay ..."
#    into just the new bold heading text:
#       "Play Crazy Colt Free - Review of JVC's Western-Themed Slot"
#    keeping only the leading empty run + the bold run.
# ------------------------------------------------------------------
$countNow = $d.Paragraphs.Count
$newHeadingPara = $d.Paragraphs.Item($countNow - 1)

$paraText = $newHeadingPara.Range.Text
$colonIdx = $paraText.IndexOf(":")

$paraStart = $newHeadingPara.Range.Start
$newTitle = "Play Crazy Colt Free - Review of JVC's Western-Themed Slot"

$boldRunRange = $d.Range($paraStart, $paraStart + $colonIdx)
$boldRunRange.Text = $newTitle

# Delete the remaining (now stale) ": Immersive gameplay..." run text.
$refreshedPara = $d.Paragraphs.Item($countNow - 1)
$remainderStart = $paraStart + $newTitle.Length
$remainderEnd = $refreshedPara.Range.End - 1
if ($remainderEnd -gt $remainderStart) {
  $remainderRange = $d.Range($remainderStart, $remainderEnd)
  $remainderRange.Delete()
}

# ------------------------------------------------------------------
# 5) Replace the body text of the final paragraph (the old image-
#    prompt paragraph) with the new meta-description text, keeping
#    its existing italic run formatting.
# ------------------------------------------------------------------
$oldImagePrompt = 'Create a feature image that stands out with a cartoon-style happy Maya warrior. The image should be set against the dusty landscape of Arizona, mirroring the theme of the game, "Crazy Colt". The warrior should be wearing glasses to highlight their intelligence, adding a unique character trait to the image. Make sure to convey a sense of excitement and adventure in the image, as the slot game is all about thrilling moments and big wins. Use bold colors and sharp lines to make the image stand out, attracting attention to this game and enticing players to take a spin.'
$newMetaDescription = 'Immersive gameplay with a Western twist. Review covers symbols, paylines, RTP, and features. Play Crazy Colt free.'

$d.Content.Find.Execute($oldImagePrompt, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newMetaDescription, 2)
